$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Season" column header
$ws.Range("AG1").Value = "Season"

# Fill the Season value (2023) for every data row
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 33).Value = 2023
}

# Match the author's final on-screen selection over the new column
$ws.Range("AG2:AG30").Select() | Out-Null
